$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price (column D) updates ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "265.38"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.72"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.283"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06157"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.593"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.718"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.350"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8290"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1589"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08292"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03394"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03138"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09235"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.896"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001711"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04771"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006209"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.005940"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04639"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006953"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1137"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003400"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01073"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006143"

# --- Rows 21-27: coin list rotated by one position with refreshed price/volume data ---
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.001089"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001500"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.769"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.292"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3380"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"
$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1234"
$ws.Range("E26").Value = "25ProBitTokenPROB"
$ws.Range("B27").Value = "UpBots"
$ws.Range("C27").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002681"
$ws.Range("E27").Value = "26UpBotsUBXT"

# --- Row 48: BOLO price + volume label update ---
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1989"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
